# Apply "curvature info for incourse" edit to 区間情報インコースまとめ.xlsx
#
# Summary of the underlying authorial change (from the OOXML diff):
#  - Sheet "数値データ": several "null" placeholder cells in column D (rows 9,
#    11, 12, 13) are replaced by real curvature-derived distance formulas
#    (D = I * (PI()/180) * G), and their supporting radius-of-curvature
#    inputs (G/H/I columns) are updated with new measured values.
#  - E13 (a stray "null") is cleared out entirely.
#  - E14 gets a new annotation value "638?".
#  - The active selection on "数値データ" moves to D15.
#  - On sheet "コース図" the view is scrolled down and zoomed in slightly,
#    and a small red line-connector shape is nudged (moved down/right and
#    enlarged a touch).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("数値データ")
$ws2 = $wb.Worksheets.Item("コース図")

# ---------------------------------------------------------------------
# Sheet "数値データ" - curvature formulas
# ---------------------------------------------------------------------
$ws1.Activate()

# Row 9 (FirstCurve)
$ws1.Range("D9").Formula = "=I9*(PI()/180)*G9"
$ws1.Range("H9").Formula = "=1/760"
# I9 keeps its existing formula "=1/H9" and recalculates automatically.

# Row 11 (SecondCurve_1)
$ws1.Range("D11").Formula = "=I11*(PI()/180)*G11"
$ws1.Range("H11").Formula = "=1/485"

# Row 12 (SecondCurve_2)
$ws1.Range("D12").Formula = "=I12*PI()/180*G12"
$ws1.Range("G12").Formula = "=90-20.63"
$ws1.Range("H12").Formula = "=1/(370)"

# Row 13 (SecondCurve_3)
$ws1.Range("D13").Formula = "=I13*PI()/180*G13"
$ws1.Range("E13").ClearContents()
$ws1.Range("G13").Formula = "=27.34+24.63"
$ws1.Range("H13").Formula = "=1/330"

# Row 14 - new annotation
$ws1.Range("E14").Value = "638?"

# Move the active selection to D15, matching the edited workbook.
$ws1.Range("D15").Select()

# ---------------------------------------------------------------------
# Sheet "コース図" - view + shape nudge
# ---------------------------------------------------------------------
$ws2.Activate()

$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1

$shp = $ws2.Shapes.Item("直線コネクタ 26")
$shp.Left = 125.18527559055119
$shp.Top = 269.6412598425197
$shp.Width = 68.32047244094488
$shp.Height = 22.695669291338582

# コース図 stays the active tab (matches activeTab="1" in the saved file).
$ws2.Activate()

Write-Output "edit applied"
